$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Spring_2020): status future -> current; priority gains 0.6 ---
$ws.Range("B2").Value = "current"

# D2 was blank (no number format applied); pick up the numeric/General
# formatting already used by the rest of column D (copy format from D3,
# which carries that style) before writing the new priority value.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D2").Value = 0.6

# --- Row 3 (Fall_2019): status current -> past; lastmod gains a date; priority drops ---
$ws.Range("B3").Value = "past"

# C3 was blank. Assigning an actual date value makes Excel apply its
# built-in short-date number format (matching the date cells below it)
# instead of leaving the cell as General.
$ws.Range("C3").Value = [datetime]"2018-12-10"
$ws.Range("D3").Value = 0.3

# --- Row 4 (Spring_2019): priority drops, lastmod unchanged ---
$ws.Range("D4").Value = 0.1

# --- Row 5 (Fall_2018): priority drops, lastmod unchanged ---
$ws.Range("D5").Value = 0.05

# Selection left where the edits finished, on B6.
$ws.Range("B6").Select()
